$d = $word.ActiveDocument

$replacements = @(
    @("146×9=", "124×6="),
    @("295×7=", "996×3="),
    @("351×6=", "321×3="),
    @("544×6=", "132×7="),
    @("160×8=", "571×8="),
    @("154×2=", "357×9="),
    @("441×6=", "386×4="),
    @("443×6=", "759×6="),
    @("532×4=", "624×6="),
    @("540×8=", "744×2="),
    @("773×6=", "409×3="),
    @("585×6=", "918×6="),
    @("139×3=", "675×9="),
    @("761×8=", "356×7="),
    @("899×7=", "739×6="),
    @("163×5=", "631×4="),
    @("144×4=", "452×9="),
    @("509×6=", "402×5="),
    @("108×7=", "681×8="),
    @("511×9=", "417×3="),
    @("472×9=", "823×6="),
    @("881×3=", "242×3="),
    @("174×3=", "861×9="),
    @("678×6=", "596×8=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
